$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fill in shared values (insertion order chosen to mirror the author's
#     edit order so new shared-string indices line up with the target) ---

# "Check storyline paper" used by D10/D11/D12
$ws.Range("D10").Value() = "Check storyline paper"
$ws.Range("D11").Value() = "Check storyline paper"
$ws.Range("D12").Value() = "Check storyline paper"

# Column C (road type descriptions) for rows 4-8
$ws.Range("C4").Value() = "A restricted access major divided highway, normally with 2 or more running lanes plus emergency hard shoulder. "
$ws.Range("C5").Value() = "The most important roads in a country's system that aren't motorways."
$ws.Range("C6").Value() = "The next most important roads in a country's system, often linking larger towns."
$ws.Range("C7").Value() = "The next most important roads in a country's system, often linking towns."
$ws.Range("C8").Value() = "The next most important roads in a country's system, often linking smaller towns and villages."

# Column D (storyline references)
$ws.Range("D4").Value() = "['F7.4', 'F7.5', 'F7.6', 'F7.7']"
$ws.Range("D5").Value() = "['F7.4', 'F7.5', 'F7.6', 'F7.7']"
$ws.Range("D6").Value() = "[F7.8, F7.9]"
$ws.Range("D7").Value() = "[F7.8, F7.9]"
$ws.Range("D8").Value() = "[F7.8, F7.9]"
$ws.Range("D9").Value() = "[F7.8, F7.9]"

# Column F5 / J4
$ws.Range("F5").Value() = "Flow velocity. Abudance of sophisticated accessories. "
$ws.Range("J4").Value() = "There are more road curves. Include them as well?"

# Column E (source reference lists)
$ws.Range("E8").Value() = "[72]"
$ws.Range("E6").Value() = "[70, 75]"
$ws.Range("E5").Value() = "[69, 81]"
$ws.Range("E4").Value() = "[68, 74, 82, 89, 90]"
$ws.Range("E7").Value() = "[71, 76, 88]"

# Column F4 / E9
$ws.Range("F4").Value() = "Fragility functions (Flow velocity. Abudance of sophisticated accessories.) and construction costs"
$ws.Range("E9").Value() = "[73, 77, 83]"

# Column F (remaining rows reuse the same fragility-function note)
$ws.Range("F6").Value() = "Fragility functions (Flow velocity)."
$ws.Range("F7").Value() = "Fragility functions (Flow velocity)."
$ws.Range("F8").Value() = "Fragility functions (Flow velocity)."
$ws.Range("F9").Value() = "Fragility functions (Flow velocity)."

# --- Formatting ---

# Wrap text in column C for rows 4-8 (no hyperlink on these, matches style index 1)
$ws.Range("C4").WrapText = $true
$ws.Range("C5").WrapText = $true
$ws.Range("C6").WrapText = $true
$ws.Range("C7").WrapText = $true
$ws.Range("C8").WrapText = $true

# Rows 6-9 column D loses its previous wrap-text styling (goes back to default)
$ws.Range("D6").ClearFormats()
$ws.Range("D6").Value() = "[F7.8, F7.9]"
$ws.Range("D7").ClearFormats()
$ws.Range("D7").Value() = "[F7.8, F7.9]"
$ws.Range("D8").ClearFormats()
$ws.Range("D8").Value() = "[F7.8, F7.9]"
$ws.Range("D9").ClearFormats()
$ws.Range("D9").Value() = "[F7.8, F7.9]"

# Row heights
$ws.Rows.Item(4).RowHeight = 58
$ws.Rows.Item(5).RowHeight = 29
$ws.Rows.Item(6).RowHeight = 43.5
$ws.Rows.Item(7).RowHeight = 29
$ws.Rows.Item(8).RowHeight = 43.5
$ws.Rows.Item(9).RowHeight = 29

# --- Hyperlinks: keep only C9, C10, C11 (drop C4-C8) ---
# (Hyperlinks.Add forces the built-in "Hyperlink" style onto the cell, so
#  reset back to Normal + re-apply wrap text afterwards for these 3 cells.
#  Add(... TextToDisplay) also clobbers the cell's actual text, so restore
#  the original cell value afterwards - the hyperlink keeps its separate
#  "display" attribute regardless.)
$origC11 = $ws.Range("C11").Value()
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C10"), "\")
$ws.Hyperlinks.Add($ws.Range("C9"), "\")
$ws.Hyperlinks.Add($ws.Range("C11"), "\", $null, $null, "\\")
$ws.Range("C11").Value() = $origC11

$ws.Range("C10").Style = "Normal"
$ws.Range("C10").WrapText = $true
$ws.Range("C9").Style = "Normal"
$ws.Range("C9").WrapText = $true
$ws.Range("C11").Style = "Normal"
$ws.Range("C11").WrapText = $true

# --- Sheet view: zoom 70 -> 85, selection E14 -> G6, drop frozen topLeftCell ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 85
$ws.Range("G6").Select()
